$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 111, shifting existing rows 111-139 down to 112-140.
$ws.Rows("111").Insert()

# Populate the newly inserted row 111 with its data (mirrors the constant
# columns from the surrounding rows, plus the new record's own values).
$ws.Cells.Item(111, 1).Value = 8
$ws.Cells.Item(111, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(111, 3).Value = "Coquimbo"
$ws.Cells.Item(111, 4).Value = 44722
$ws.Cells.Item(111, 5).Value = 4
$ws.Cells.Item(111, 6).Value = 100112001
$ws.Cells.Item(111, 7).Value = "Berenjena"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 440
$ws.Cells.Item(111, 11).Value = 8000
$ws.Cells.Item(111, 12).Value = 9000
$ws.Cells.Item(111, 13).Value = 8500
$ws.Cells.Item(111, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(111, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(111, 16).Value = 170
$ws.Cells.Item(111, 17).Value = 50
$ws.Cells.Item(111, 18).Value = "Hortaliza"
